$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week date range) ---
$ws.Range("A8").Value = "Volume 30   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/11/2023  Through  9/17/2023"

# --- Cells that switch between the text "N/A" markers and real numbers ---
# (PasteSpecial formats-only from an untouched reference cell keeps the
#  correct shared style index, matching the sibling numeric/text cells.)
$ws.Range("I14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 1

$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G22").PasteSpecial(-4122)

$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H22").PasteSpecial(-4122)

$ws.Range("I14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = 1

$ws.Range("K14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = 0

$ws.Range("I14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1

$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F30").PasteSpecial(-4122)

$ws.Range("I14").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1

$ws.Range("K14").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Value = -100

$excel.CutCopyMode = 0

# --- Remaining plain numeric value updates ---
$ws.Range("I15").Value = 14
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 0
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -37.5
$ws.Range("I16").Value = 59
$ws.Range("J16").Value = 66
$ws.Range("K16").Value = -10.60606060606
$ws.Range("L16").Value = 18
$ws.Range("M16").Value = -53.174603174603
$ws.Range("N16").Value = -88.339920948616
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = -32.142857142857
$ws.Range("I17").Value = 153
$ws.Range("J17").Value = 164
$ws.Range("K17").Value = -6.70731707317
$ws.Range("L17").Value = 37.837837837837
$ws.Range("M17").Value = 34.210526315789
$ws.Range("N17").Value = -35.443037974683
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 9
$ws.Range("J18").Value = 136
$ws.Range("K18").Value = -38.970588235294
$ws.Range("L18").Value = -23.853211009174
$ws.Range("M18").Value = -69.818181818181
$ws.Range("N18").Value = -93.191140278917
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 7.692307692307
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = 34.210526315789
$ws.Range("I19").Value = 410
$ws.Range("J19").Value = 441
$ws.Range("K19").Value = -7.029478458049
$ws.Range("L19").Value = 34.868421052631
$ws.Range("M19").Value = 46.953405017921
$ws.Range("N19").Value = -12.20556745182
$ws.Range("C20").Value = 5
$ws.Range("E20").Value = 25
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 111
$ws.Range("J20").Value = 82
$ws.Range("K20").Value = 35.365853658536
$ws.Range("L20").Value = 88.135593220339
$ws.Range("M20").Value = 7.766990291262
$ws.Range("N20").Value = -91.36186770428
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -23.529411764705
$ws.Range("F21").Value = 98
$ws.Range("G21").Value = 108
$ws.Range("H21").Value = -9.259259259259
$ws.Range("I21").Value = 832
$ws.Range("J21").Value = 904
$ws.Range("K21").Value = -7.964601769911
$ws.Range("L21").Value = 29.797191887675
$ws.Range("M21").Value = -7.964601769911
$ws.Range("N21").Value = -77.730192719486
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = -26.923076923076
$ws.Range("F24").Value = 71
$ws.Range("G24").Value = 112
$ws.Range("H24").Value = -36.607142857142
$ws.Range("I24").Value = 781
$ws.Range("J24").Value = 811
$ws.Range("K24").Value = -3.699136868064
$ws.Range("L24").Value = 36.538461538461
$ws.Range("M24").Value = 19.236641221374
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -73.333333333333
$ws.Range("F25").Value = 25
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = -45.652173913043
$ws.Range("I25").Value = 247
$ws.Range("J25").Value = 250
$ws.Range("K25").Value = -1.2
$ws.Range("L25").Value = 15.962441314554
$ws.Range("M25").Value = -19.016393442623
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 16
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = -11.111111111111
$ws.Range("L26").Value = 45.454545454545
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -40
$ws.Range("J27").Value = 61
$ws.Range("K27").Value = -19.672131147541
$ws.Range("L27").Value = -5.76923076923
$ws.Range("N28").Value = -75
$ws.Range("N29").Value = -75
$ws.Range("E30").Value = -100
$ws.Range("J30").Value = 12
$ws.Range("K30").Value = -50
